# Apply crypto price/volume update (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.299.26"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.238.82"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.20"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.641"
$ws.Range("E6").Value = "  +2.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.72"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.438"
$ws.Range("E9").Value = "  +3.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0950"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.38"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.97"
$ws.Range("E12").Value = "  +5.64%  "
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.576.14"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.18"
$ws.Range("E15").Value = "  -2.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.00"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.821"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.250.76"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.293.69"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0961"
$ws.Range("E20").Value = "  -4.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.92"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.74"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.76"
$ws.Range("E25").Value = "  +33.53%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.29"
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "174.27"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.62"
$ws.Range("E30").Value = "  +5.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.129"
$ws.Range("E31").Value = "  -5.14%  "
$ws.Range("E32").Value = "  +2.19%  "
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("E34").Value = "  +4.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0673"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.88"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.57"
$ws.Range("E37").Value = "  -6.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.30"
$ws.Range("E38").Value = "  -5.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0249"
$ws.Range("E40").Value = "  +0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.58"
$ws.Range("E42").Value = "  +4.05%  "
$ws.Range("E43").Value = "  +4.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.95"
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.27"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0938"
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000207"
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.427.96"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("B50").Value = "Celestia"
$ws.Range("C50").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.84"
$ws.Range("E50").Value = "  +4.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.25"
$ws.Range("E51").Value = "  -1.51%  "
